# Insert a new data row at row 623, shifting the existing row 623 (and all
# rows below it, through row 708) down by one position to row 624..709.
# The values below for A623:R623 are the brand new record added to the
# dataset; everything previously at row 623 onward is preserved, just moved
# down one row (Excel's native row-insert behavior handles that shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 623, pushing old rows 623-708
# down to 624-709. Excel carries the row-above's formatting into the new
# row automatically, so the D (Fecha) column keeps its date number format.
$ws.Rows("623:623").Insert()

$ws.Cells.Item(623, 1).Value  = 6
$ws.Cells.Item(623, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(623, 3).Value  = "Metropolitana"
$ws.Cells.Item(623, 4).Value  = 44984
$ws.Cells.Item(623, 5).Value  = 13
$ws.Cells.Item(623, 6).Value  = 100112012
$ws.Cells.Item(623, 7).Value  = "Espinaca"
$ws.Cells.Item(623, 8).Value  = "Sin especificar"
$ws.Cells.Item(623, 9).Value  = "Primera"
$ws.Cells.Item(623, 10).Value = 370
$ws.Cells.Item(623, 11).Value = 7000
$ws.Cells.Item(623, 12).Value = 7500
$ws.Cells.Item(623, 13).Value = 7297
$ws.Cells.Item(623, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(623, 15).Value = "Región Metropolitana"
$ws.Cells.Item(623, 16).Value = 730
$ws.Cells.Item(623, 17).Value = 10
$ws.Cells.Item(623, 18).Value = "Hortaliza"
